# Exploring Self-Management.pptx
# "Double check the Key card and Fixes #4"
#
# 1) Scrum Master "Key card": capitalize "scrum" -> "Scrum" in
#    "Facilitates the scrum events" (appears on two duplicated slides).
# 2) Product Owner "Key card": shorten
#    "Focus on creating value & continuously validate assumptions"
#    to "Focus on creating value & process" (appears on two duplicated
#    slides; on the second one the replacement text ends up as its own
#    run because it was typed over the selected tail of the sentence).

$p = $ppt.ActivePresentation

# --- Slide 17: Scrum Master key card ---
$p.Slides.Item(17).Shapes.Item(1).TextFrame.TextRange.Text = "Facilitates the Scrum events"

# --- Slide 18: duplicate Scrum Master key card ---
$p.Slides.Item(18).Shapes.Item(1).TextFrame.TextRange.Text = "Facilitates the Scrum events"

# --- Slide 31: Product Owner key card ---
$p.Slides.Item(31).Shapes.Item(1).TextFrame.TextRange.Text = "Focus on creating value & process"

# --- Slide 32: duplicate Product Owner key card ---
# Keep the leading "Focus on creating value " run untouched and retype
# just the trailing clause, so it lands in its own run.
$tr32 = $p.Slides.Item(32).Shapes.Item(1).TextFrame.TextRange
$tail = $tr32.Characters(25, 35)
$tail.Text = "& process"
